$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 300
$ws.Range("I13").Value = 300
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -131

$ws.Range("H34").Value = 2430.25
$ws.Range("I34").Value = 2430.25
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2430.25
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2227.25

$ws.Range("H36").Value = 2430.25
$ws.Range("I36").Value = 2430.25
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2430.25
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1715.25

$ws.Range("H132").Value = 21514.4
$ws.Range("I132").Value = 21514.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 64543.2
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -62013.2
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1910.0741
$ws.Range("I22").Value = 1639.4667
$ws.Range("J22").Value = 2248.3333
$ws.Range("K22").Value = 1639.4667
$ws.Range("L22").Value = 2248.3333
$ws.Range("M22").Value = -1340.4667
$ws.Range("N22").Value = -2846.3333

$ws.Range("H112").Value = 45000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 45000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 45000
$ws.Range("N112").Value = -47954

$ws.Range("H122").Value = 2636.6667
$ws.Range("I122").Value = 1605
$ws.Range("J122").Value = 4700
$ws.Range("K122").Value = 4815
$ws.Range("L122").Value = 14100
$ws.Range("M122").Value = -2365
$ws.Range("N122").Value = -19000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1899.6666
$ws.Range("I5").Value = 349.5
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 349.5
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -236.5
$ws.Range("N5").Value = -5226

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H36").Value = 2805.6
$ws.Range("I36").Value = 2805.6
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2805.6
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2271.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 5001500
$ws.Range("I3").Value = 5001500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5001500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5001387

$ws.Range("H22").Value = 480
$ws.Range("I22").Value = 480
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 480
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -130

$ws.Range("H31").Value = 5810.891
$ws.Range("I31").Value = 3033
$ws.Range("J31").Value = 6075.452
$ws.Range("K31").Value = 3033
$ws.Range("L31").Value = 6075.452
$ws.Range("M31").Value = -2738
$ws.Range("N31").Value = -6665.452

$ws.Range("H34").Value = 5810.891
$ws.Range("I34").Value = 3033
$ws.Range("J34").Value = 6075.452
$ws.Range("K34").Value = 3033
$ws.Range("L34").Value = 6075.452
$ws.Range("M34").Value = -2831
$ws.Range("N34").Value = -6479.452

$ws.Range("H94").Value = 4171.231
$ws.Range("I94").Value = 2274.625
$ws.Range("J94").Value = 7205.8
$ws.Range("K94").Value = 2274.625
$ws.Range("L94").Value = 7205.8
$ws.Range("M94").Value = -1823.625
$ws.Range("N94").Value = -8107.8

$ws.Range("H122").Value = 1268.5625
$ws.Range("I122").Value = 1175.4166
$ws.Range("J122").Value = 1548
$ws.Range("K122").Value = 3526.2498
$ws.Range("L122").Value = 4644
$ws.Range("M122").Value = -1076.2498
$ws.Range("N122").Value = -9544

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 6915.4443
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 6915.4443
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 20746.3329
$ws.Range("N17").Value = -21084.3329
$ws.Range("M17").ClearContents()

$ws.Range("H39").Value = 6382.0713
$ws.Range("I39").Value = 315
$ws.Range("J39").Value = 7393.25
$ws.Range("K39").Value = 945
$ws.Range("L39").Value = 22179.75
$ws.Range("M39").Value = -651
$ws.Range("N39").Value = -22767.75

$ws.Range("H58").Value = 2000
$ws.Range("I58").Value = 2000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 6000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -5872

$ws.Range("H80").Value = 4169.391
$ws.Range("I80").Value = 3784.0527
$ws.Range("J80").Value = 5999.75
$ws.Range("K80").Value = 11352.1581
$ws.Range("L80").Value = 17999.25
$ws.Range("M80").Value = -10416.1581
$ws.Range("N80").Value = -19871.25

$ws.Range("H81").Value = 1469.8
$ws.Range("I81").Value = 950
$ws.Range("J81").Value = 2249.5
$ws.Range("K81").Value = 2850
$ws.Range("L81").Value = 6748.5
$ws.Range("M81").Value = -1727
$ws.Range("N81").Value = -8994.5

$ws.Range("H83").Value = 4169.391
$ws.Range("I83").Value = 3784.0527
$ws.Range("J83").Value = 5999.75
$ws.Range("K83").Value = 34056.4743
$ws.Range("L83").Value = 53997.75
$ws.Range("M83").Value = -29376.4743
$ws.Range("N83").Value = -63357.75

$ws.Range("H84").Value = 1469.8
$ws.Range("I84").Value = 950
$ws.Range("J84").Value = 2249.5
$ws.Range("K84").Value = 8550
$ws.Range("L84").Value = 20245.5
$ws.Range("M84").Value = -2934
$ws.Range("N84").Value = -31477.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 47221.75
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 47221.75
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 47221.75
$ws.Range("N20").Value = -47711.75

$ws.Range("H24").Value = 46666.2
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 46666.2
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 46666.2
$ws.Range("N24").Value = -47012.2
$ws.Range("M24").ClearContents()

$ws.Range("H102").Value = 2181.7693
$ws.Range("I102").Value = 2011.3
$ws.Range("J102").Value = 2750
$ws.Range("K102").Value = 2011.3
$ws.Range("L102").Value = 2750
$ws.Range("M102").Value = -389.3
$ws.Range("N102").Value = -5994

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 29999
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 29999
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 29999
$ws.Range("N39").Value = -30919

$ws.Range("H61").Value = 3388.5454
$ws.Range("I61").Value = 1659.25
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 1659.25
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -1457.25
$ws.Range("N61").Value = -8404

$ws.Range("H68").Value = 7666
$ws.Range("I68").Value = 2998
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 2998
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -2249
$ws.Range("N68").Value = -11498

$ws.Range("H71").Value = 7666
$ws.Range("I71").Value = 2998
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 14990
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -11246
$ws.Range("N71").Value = -57488

$ws.Range("H87").Value = 40000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 40000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42246

$ws.Range("H90").Value = 40000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 40000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -131232

$ws.Range("H110").Value = 45000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 45000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 45000
$ws.Range("N110").Value = -53180

$ws.Range("H113").Value = 3388.5454
$ws.Range("I113").Value = 1659.25
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 1659.25
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = 510.75
$ws.Range("N113").Value = -12340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 99908
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 99908
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 99908
$ws.Range("N80").Value = -101904

$ws.Range("H83").Value = 99908
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 99908
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 299724
$ws.Range("N83").Value = -309708

$ws.Range("H100").Value = 660.2143
$ws.Range("I100").Value = 594.3
$ws.Range("J100").Value = 825
$ws.Range("K100").Value = 1188.6
$ws.Range("L100").Value = 1650
$ws.Range("M100").Value = -647.5999999999999
$ws.Range("N100").Value = -2732

$ws.Range("H122").Value = 5479
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 5598.75
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 16796.25
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -21696.25
